# Fruta / hortaliza, semanal
# Insert a new weekly price-survey row for "Brócoli" at Feria Lagunitas de
# Puerto Montt, pushing the existing rows 316-333 down to 317-334.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 316 (shifts rows 316:333 -> 317:334,
# and grows the sheet's used range to row 334).
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A316").Value = 4
$ws.Range("B316").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C316").Value = "Los Lagos"
$ws.Range("D316").Value = 44706
$ws.Range("E316").Value = 10
$ws.Range("F316").Value = 100112023
$ws.Range("G316").Value = "Brócoli"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 250
$ws.Range("K316").Value = 1500
$ws.Range("L316").Value = 1500
$ws.Range("M316").Value = 1500
$ws.Range("N316").Value = "$/unidad"
$ws.Range("O316").Value = "Región del Maule"
$ws.Range("P316").Value = 1500
$ws.Range("Q316").Value = 1
$ws.Range("R316").Value = "Hortaliza"
